# Updates the cryptocurrency price/volume table on Sheet1 to the latest
# scraped snapshot (GitHub Actions data refresh).
#
# Column D ("Price") and E ("Volume(1h)") hold text-formatted values (some
# look numeric, e.g. "220.99"), so a leading apostrophe forces Excel to
# keep them as text instead of silently coercing to a Double (which would
# corrupt values like "1.00" -> 1 or "0.170" -> 0.17000000000000001).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.454.20'
$ws.Range("E2").Value = '  +4.18%  '
$ws.Range("D3").Value = '3.130.23'
$ws.Range("E3").Value = '  +2.44%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'220.99"
$ws.Range("E5").Value = '  +6.12%  '
$ws.Range("D6").Value = "'623.85"
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = "'0.974"
$ws.Range("E7").Value = '  +21.60%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = "'0.378"
$ws.Range("E8").Value = '  +3.83%  '
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '3.126.52'
$ws.Range("E10").Value = '  +2.41%  '
$ws.Range("D11").Value = "'0.745"
$ws.Range("E11").Value = '  +26.85%  '
$ws.Range("E12").Value = '  +6.73%  '
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = '  +8.57%  '
$ws.Range("D14").Value = "'35.01"
$ws.Range("E14").Value = '  +11.53%  '
$ws.Range("E15").Value = '  +4.01%  '
$ws.Range("D16").Value = '91.124.24'
$ws.Range("E16").Value = '  +3.91%  '
$ws.Range("D17").Value = '3.707.85'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").Value = '3.119.48'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("D19").Value = "'3.77"
$ws.Range("E19").Value = '  +19.02%  '
$ws.Range("D20").Value = "'0.0000227"
$ws.Range("E20").Value = '  +17.08%  '
$ws.Range("D21").Value = "'14.18"
$ws.Range("E21").Value = '  +8.92%  '
$ws.Range("D22").Value = "'436.92"
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("D23").Value = "'8.79"
$ws.Range("E23").Value = '  +9.41%  '
$ws.Range("E24").Value = '  +7.61%  '
$ws.Range("D25").Value = "'6.06"
$ws.Range("E25").Value = '  +13.51%  '
$ws.Range("D26").Value = "'12.38"
$ws.Range("E26").Value = '  +7.62%  '
$ws.Range("D27").Value = "'85.93"
$ws.Range("E27").Value = '  +6.61%  '
$ws.Range("D28").Value = '3.294.45'
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = "'0.170"
$ws.Range("E30").Value = '  +7.26%  '
$ws.Range("D31").Value = "'9.02"
$ws.Range("E31").Value = '  +13.18%  '
$ws.Range("D32").Value = "'530.80"
$ws.Range("E32").Value = '  +6.57%  '
$ws.Range("D33").Value = "'3.84"
$ws.Range("E33").Value = '  +10.19%  '
$ws.Range("D34").Value = "'0.896"
$ws.Range("E34").Value = '  -17.57%  '
$ws.Range("D35").Value = "'7.22"
$ws.Range("E35").Value = '  +10.83%  '
$ws.Range("D36").Value = "'0.146"
$ws.Range("E36").Value = '  +12.08%  '
$ws.Range("D37").Value = "'23.79"
$ws.Range("E37").Value = '  +8.32%  '
$ws.Range("D38").Value = "'1.29"
$ws.Range("E38").Value = '  +6.17%  '
$ws.Range("D39").Value = "'1.87"
$ws.Range("E39").Value = '  +5.01%  '
$ws.Range("D40").Value = "'22.31"
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").Value = "'0.155"
$ws.Range("E42").Value = '  +17.85%  '
$ws.Range("D43").Value = "'0.0791"
$ws.Range("E43").Value = '  +18.46%  '
$ws.Range("D44").Value = "'0.384"
$ws.Range("E44").Value = '  +8.17%  '
$ws.Range("E46").Value = '  +8.87%  '
$ws.Range("D47").Value = "'145.92"
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = "'44.12"
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("E49").Value = '  +11.81%  '
$ws.Range("D50").Value = "'168.44"
$ws.Range("E50").Value = '  +9.08%  '
$ws.Range("E51").Value = '  +25.20%  '
